# Update "想去人数" (F column) values on the "展览" (sheet1) and
# "全部类型" (sheet4) worksheets, per regenerated site data (456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value  = 215
$wsExhibition.Range("F7").Value  = 116
$wsExhibition.Range("F10").Value = 44
$wsExhibition.Range("F11").Value = 7007
$wsExhibition.Range("F14").Value = 3440
$wsExhibition.Range("F15").Value = 252
$wsExhibition.Range("F16").Value = 445
$wsExhibition.Range("F18").Value = 581
$wsExhibition.Range("F19").Value = 59

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 215
$wsAll.Range("F9").Value  = 116
$wsAll.Range("F12").Value = 44
$wsAll.Range("F14").Value = 7007
$wsAll.Range("F18").Value = 3440
$wsAll.Range("F19").Value = 252
$wsAll.Range("F20").Value = 445
$wsAll.Range("F22").Value = 581
$wsAll.Range("F23").Value = 59
